$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1 (2)")

# Row 21: E21 gets a value of -5 (was blank)
$ws.Range("E21").Value = -5

# Row 55: clear C55 and D55 (they become blank cells)
$ws.Range("C55:D55").ClearContents()

# Row 68: G68 text changes from "F72A" equivalent old shared string (index 162) to new string "M323G" (index 223, new unique string)
$ws.Range("G68").Value = "M323G"

# Row 94: B94 value changes 0.13 -> 0.1; C94:E94 cleared
$ws.Range("B94").Value = 0.1
$ws.Range("C94:E94").ClearContents()

# Selection change
$ws.Range("AD8").Select()
